$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = [double]"1.005539333333333"
$ws.Range("H2").Value = [double]"3.016618"
$ws.Range("I2").Value = [double]"0.06022996424239081"
$ws.Range("J2").Value = [double]"0.06022996424239079"
$ws.Range("K2").Value = [double]"3"
$ws.Range("L2").Value = [double]"1"
$ws.Range("M2").Value = [double]"0.1178306666666667"
$ws.Range("N2").Value = [double]"0.353492"
$ws.Range("O2").Value = [double]"0.0005211703885903252"
$ws.Range("P2").Value = [double]"0.0005211703885903251"
$ws.Range("Q2").Value = [double]"0.1184833700062222"
$ws.Range("R2").Value = [double]"1.066350330056"
$ws.Range("S2").Value = [double]"3.139007386898821E-05"
$ws.Range("T2").Value = [double]"3.139007386898819E-05"
$ws.Range("G3").Value = [double]"1.005539333333333"
$ws.Range("H3").Value = [double]"3.016618"
$ws.Range("I3").Value = [double]"0.06022996424239081"
$ws.Range("J3").Value = [double]"0.06022996424239079"
$ws.Range("O3").Value = [double]"0.9986266812609277"
$ws.Range("P3").Value = [double]"0.9986266812609277"
$ws.Range("Q3").Value = [double]"227.0287360223226"
$ws.Range("R3").Value = [double]"2043.258624200904"
$ws.Range("S3").Value = [double]"0.06014724930384308"
$ws.Range("T3").Value = [double]"0.06014724930384306"
$ws.Range("G4").Value = [double]"1.005539333333333"
$ws.Range("H4").Value = [double]"3.016618"
$ws.Range("I4").Value = [double]"0.06022996424239081"
$ws.Range("J4").Value = [double]"0.06022996424239079"
$ws.Range("M4").Value = [double]"0.192661"
$ws.Range("N4").Value = [double]"0.5779829999999999"
$ws.Range("O4").Value = [double]"0.0008521483504820529"
$ws.Range("P4").Value = [double]"0.0008521483504820528"
$ws.Range("Q4").Value = [double]"0.1937282134993333"
$ws.Range("R4").Value = [double]"1.743553921494"
$ws.Range("S4").Value = [double]"5.132486467874635E-05"
$ws.Range("T4").Value = [double]"5.132486467874634E-05"
$ws.Range("H5").Value = [double]"4.211987000000001"
$ws.Range("I5").Value = [double]"0.08409676876535742"
$ws.Range("J5").Value = [double]"0.0840967687653574"
$ws.Range("K5").Value = [double]"3"
$ws.Range("L5").Value = [double]"1"
$ws.Range("M5").Value = [double]"0.1178306666666667"
$ws.Range("N5").Value = [double]"0.353492"
$ws.Range("O5").Value = [double]"0.0005211703885903252"
$ws.Range("P5").Value = [double]"0.0005211703885903251"
$ws.Range("Q5").Value = [double]"0.1654337454004445"
$ws.Range("R5").Value = [double]"1.488903708604"
$ws.Range("S5").Value = [double]"4.382874565663205E-05"
$ws.Range("T5").Value = [double]"4.382874565663204E-05"
$ws.Range("H6").Value = [double]"4.211987000000001"
$ws.Range("I6").Value = [double]"0.08409676876535742"
$ws.Range("J6").Value = [double]"0.0840967687653574"
$ws.Range("O6").Value = [double]"0.9986266812609277"
$ws.Range("P6").Value = [double]"0.9986266812609277"
$ws.Range("S6").Value = [double]"0.08398127709691652"
$ws.Range("T6").Value = [double]"0.08398127709691651"
$ws.Range("H7").Value = [double]"4.211987000000001"
$ws.Range("I7").Value = [double]"0.08409676876535742"
$ws.Range("J7").Value = [double]"0.0840967687653574"
$ws.Range("M7").Value = [double]"0.192661"
$ws.Range("N7").Value = [double]"0.5779829999999999"
$ws.Range("O7").Value = [double]"0.0008521483504820529"
$ws.Range("P7").Value = [double]"0.0008521483504820528"
$ws.Range("Q7").Value = [double]"0.2704952091356667"
$ws.Range("R7").Value = [double]"2.434456882221"
$ws.Range("S7").Value = [double]"7.166292278426996E-05"
$ws.Range("T7").Value = [double]"7.166292278426994E-05"
$ws.Range("G8").Value = [double]"14.28546633333333"
$ws.Range("H8").Value = [double]"42.856399"
$ws.Range("I8").Value = [double]"0.8556732669922519"
$ws.Range("J8").Value = [double]"0.8556732669922518"
$ws.Range("K8").Value = [double]"3"
$ws.Range("L8").Value = [double]"1"
$ws.Range("M8").Value = [double]"0.1178306666666667"
$ws.Range("N8").Value = [double]"0.353492"
$ws.Range("O8").Value = [double]"0.0005211703885903252"
$ws.Range("P8").Value = [double]"0.0005211703885903251"
$ws.Range("Q8").Value = [double]"1.683266021700889"
$ws.Range("R8").Value = [double]"15.149394195308"
$ws.Range("S8").Value = [double]"0.000445951569064705"
$ws.Range("T8").Value = [double]"0.0004459515690647048"
$ws.Range("G9").Value = [double]"14.28546633333333"
$ws.Range("H9").Value = [double]"42.856399"
$ws.Range("I9").Value = [double]"0.8556732669922519"
$ws.Range("J9").Value = [double]"0.8556732669922518"
$ws.Range("O9").Value = [double]"0.9986266812609277"
$ws.Range("P9").Value = [double]"0.9986266812609277"
$ws.Range("Q9").Value = [double]"3225.345103502774"
$ws.Range("R9").Value = [double]"29028.10593152497"
$ws.Range("S9").Value = [double]"0.8544981548601682"
$ws.Range("T9").Value = [double]"0.8544981548601681"
$ws.Range("G10").Value = [double]"14.28546633333333"
$ws.Range("H10").Value = [double]"42.856399"
$ws.Range("I10").Value = [double]"0.8556732669922519"
$ws.Range("J10").Value = [double]"0.8556732669922518"
$ws.Range("M10").Value = [double]"0.192661"
$ws.Range("N10").Value = [double]"0.5779829999999999"
$ws.Range("O10").Value = [double]"0.0008521483504820529"
$ws.Range("P10").Value = [double]"0.0008521483504820528"
$ws.Range("Q10").Value = [double]"2.752252229246333"
$ws.Range("R10").Value = [double]"24.77027006321699"
$ws.Range("S10").Value = [double]"0.0007291605630190367"
$ws.Range("T10").Value = [double]"0.0007291605630190365"
